# Prepend "Design: " to the answer texts in the feedback table.
$d = $word.ActiveDocument

$targets = @(
    "Yes",
    "No problems",
    "Budget was good, yard was refunded.",
    "First time in 20 years automation design made all drawings, not outsourced. ",
    "Machinery design have skilled designers.",
    "Jira and ERM  doesn't work as I expected. ACAD without electrical symbols.",
    "More teamwork before purchase of different sister systems. "
)

foreach ($t in $targets) {
    $range = $d.Content
    $found = $range.Find.Execute($t, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Collapse(1)
        $range.InsertBefore("Design: ")
    }
}
